# Insert a new weekly price record for "Coliflor" at Macroferia Regional de
# Talca, pushing the existing rows 100-143 down to 101-144 (row 144 keeps
# the data that used to live in row 143).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(100).Insert()

$ws.Range("A100").Value = 5
$ws.Range("B100").Value = "Macroferia Regional de Talca"
$ws.Range("C100").Value = "Maule"
$ws.Range("D100").Value = 44466
$ws.Range("E100").Value = 7
$ws.Range("F100").Value = 100112008
$ws.Range("G100").Value = "Coliflor"
$ws.Range("H100").Value = "Sin especificar"
$ws.Range("I100").Value = "Primera"
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 600
$ws.Range("L100").Value = 600
$ws.Range("M100").Value = 600
$ws.Range("N100").Value = "$/unidad"
$ws.Range("O100").Value = "Región del Maule"
$ws.Range("P100").Value = 600
$ws.Range("Q100").Value = 1
$ws.Range("R100").Value = "Hortaliza"
